# added test with one PV-bus
# Inserts a new data row (row 43) above the existing blank separator /
# summary rows, giving the workbook one additional "start"/"end" entry
# (09:15 -> 12:00 on 2014-03-03). Excel shifts the separator + the three
# summary rows down by one and widens the SUM() range automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the blank separator row (and the summary rows below it) down by
# inserting a fresh row at 43; row 43 becomes available for the new entry.
$ws.Rows.Item(43).Insert()

$ws.Range("A43").Value = 2014
$ws.Range("B43").Value = 3
$ws.Range("C43").Value = 3
$ws.Range("D43").Value = 0.38541666666666669   # 09:15
$ws.Range("E43").Value = 0.5                   # 12:00
$ws.Range("F43").Formula = "=(E43-D43)*24*60"
$ws.Range("G43").Formula = "=F43/60"

# Matches the author's final selection in the workbook.
$ws.Range("D44").Select() | Out-Null
